$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in the workbook even though they
# look numeric, so mark the cells as Text format before assigning the new
# string value. This mirrors how a user would type over a Text-formatted cell
# in Excel and keeps the values as literal strings (preserving formats like
# trailing zeros, e.g. "243.15") instead of Excel re-interpreting them as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.424"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.569"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8111"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9356"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1420"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07439"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03270"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03055"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09332"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.874"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001584"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04673"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005953"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005904"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001262"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004905"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009508"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.142"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03955"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006191"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002903"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009278"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005190"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7504"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002290"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
